# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Row -> New Value pairs for the "展览" sheet
$sheet1Updates = @{
    2  = 62
    3  = 591
    6  = 11
    7  = 14800
    10 = 15197
    11 = 32
    12 = 8675
    13 = 327
    20 = 14
    21 = 3
    23 = 6
    28 = 54
    32 = 26
    36 = 421
    37 = 109
    38 = 5339
}

# Row -> New Value pairs for the "全部类型" sheet
$sheet4Updates = @{
    2  = 62
    3  = 591
    6  = 11
    7  = 14800
    10 = 15197
    11 = 32
    12 = 8675
    13 = 327
    21 = 14
    22 = 3
    24 = 6
    29 = 54
    35 = 26
    39 = 421
    40 = 109
    41 = 5339
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
